# Update the "nr" (column D) scores that were previously blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 98
$ws.Range("D27").Value = 89
$ws.Range("D40").Value = 89

# Scroll the view down and move the active selection, matching the
# author's last on-screen position when the file was saved.
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("C42").Select()
